$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text (string) representation instead of
# being re-interpreted as numbers/percentages by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "293.07"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.21%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "30.78"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "7.28%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.147"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.31%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07126"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "6.98%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.527"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.86%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.623"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "6.37%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.408"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.64%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9172"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.91%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1631"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.31%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07684"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "19.74%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07773"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.24%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02948"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.11%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.08996"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.15%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001594"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.04%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006562"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.33%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006196"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.96%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.482"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.02%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.242"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.34%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.11%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1363"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "5.05%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.148"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.70%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1600"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "3.01%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04541"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.21%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001210"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.95%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004239"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "2.39%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001169"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-6.41%"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001688"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "4.34%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04412"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.56%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007045"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "4.97%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1274"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.11%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002209"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "11.61%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01319"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "8.06%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005863"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "4.38%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.734"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-11.99%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01299"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.65%"
